$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- Paragraph 4: ">>>  your stuff after this line >>>" -> split into 3 runs w/ proofErr marks
$p4 = $d.Paragraphs.Item(4)
$rng4 = $d.Range($p4.Range.Start, $p4.Range.End)
$xml4 = "<w:p $wns>" +
        "<w:r><w:t>&gt;&gt;</w:t></w:r>" +
        "<w:proofErr w:type='gramStart'/>" +
        "<w:r><w:t>&gt;  your</w:t></w:r>" +
        "<w:proofErr w:type='gramEnd'/>" +
        "<w:r><w:t xml:space='preserve'> stuff after this line &gt;&gt;&gt;</w:t></w:r>" +
        "</w:p>"
$rng4.InsertXML($xml4)

# --- Paragraph 5: "A cool quote by Edsger Dijkstra:" -> split into 3 runs w/ proofErr marks
$p5 = $d.Paragraphs.Item(5)
$rng5 = $d.Range($p5.Range.Start, $p5.Range.End)
$xml5 = "<w:p $wns>" +
        "<w:r><w:t xml:space='preserve'>A cool quote by </w:t></w:r>" +
        "<w:proofErr w:type='spellStart'/>" +
        "<w:r><w:t>Edsger</w:t></w:r>" +
        "<w:proofErr w:type='spellEnd'/>" +
        "<w:r><w:t xml:space='preserve'> Dijkstra:</w:t></w:r>" +
        "</w:p>"
$rng5.InsertXML($xml5)

# --- Paragraph 6: quote paragraph restructure -------------------------------
# Target shape:
#   <w:p><w:pPr><w:rPr>[Arial italic quote-style]</w:rPr></w:pPr>
#        <w:r><w:t>&#8220;</w:t></w:r>
#        <w:r><w:rPr>[Arial italic quote-style]</w:rPr><w:t>Computer science ... telescopes.&#8221;</w:t></w:r>
#   </w:p>
# InsertXML on a whole-paragraph range cannot create/alter a <w:pPr> (the
# paragraph-mark formatting is only settable through the real Font/
# ParagraphFormat object model in this host), so first "seed" the mark's
# run-properties with the native Font object (applies to the whole range,
# incl. the mark), then overwrite the run content precisely with InsertXML.
$rPrXml = "<w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:i/><w:iCs/><w:color w:val='4D5156'/><w:sz w:val='21'/><w:szCs w:val='21'/><w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/></w:rPr>"

$p6 = $d.Paragraphs.Item(6)
$f6 = $p6.Range.Font
$f6.Name = "Arial"
$f6.NameBi = "Arial"
$f6.Italic = $true
$f6.ItalicBi = $true
$f6.Color = 5656909
$f6.Size = 10.5
$f6.SizeBi = 10.5

$p6b = $d.Paragraphs.Item(6)
$rng6 = $d.Range($p6b.Range.Start, $p6b.Range.End)
$xml6 = "<w:p $wns>" +
        "<w:r><w:t>&#8220;</w:t></w:r>" +
        "<w:r>$rPrXml<w:t>Computer science is no more about computers than astronomy is about telescopes.&#8221;</w:t></w:r>" +
        "</w:p>"
$rng6.InsertXML($xml6)

# --- Add new paragraph 7 with second quote, after paragraph 6
$p6c = $d.Paragraphs.Item(6)
$rng6c = $d.Range($p6c.Range.Start, $p6c.Range.End)
$rng6c.InsertParagraphAfter()
$p7 = $d.Paragraphs.Item(7)
$rng7 = $d.Range($p7.Range.Start, $p7.Range.End)
$xml7 = "<w:p $wns>" +
        "<w:r>$rPrXml<w:t>&#8220;</w:t></w:r>" +
        "<w:r>$rPrXml<w:t>Everything is possible, we just need to believe in it&#8221;</w:t></w:r>" +
        "</w:p>"
$rng7.InsertXML($xml7)

# --- Add trailing empty paragraph (8th)
$p7b = $d.Paragraphs.Item(7)
$rng7b = $d.Range($p7b.Range.Start, $p7b.Range.End)
$rng7b.InsertParagraphAfter()
$p8 = $d.Paragraphs.Item(8)
$rng8 = $d.Range($p8.Range.Start, $p8.Range.End)
$rng8.InsertXML("<w:p $wns/>")
